$wb = $excel.ActiveWorkbook

# ALC row 41
$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H41").Value = 840.26666
$ws_ALC.Range("I41").Value = 625.5
$ws_ALC.Range("J41").Value = 1269.8
$ws_ALC.Range("K41").Value = 625.5
$ws_ALC.Range("L41").Value = 1269.8
$ws_ALC.Range("M41").Value = -185.5
$ws_ALC.Range("N41").Value = -2149.8

# ALC row 55
$ws_ALC.Range("H55").Value = 107
$ws_ALC.Range("I55").Value = 55
$ws_ALC.Range("J55").Value = 159
$ws_ALC.Range("K55").Value = 55
$ws_ALC.Range("L55").Value = 159
$ws_ALC.Range("M55").Value = 159
$ws_ALC.Range("N55").Value = -587

# ALC row 76
$ws_ALC.Range("H76").Value = 71431890
$ws_ALC.Range("I76").Value = 100003040
$ws_ALC.Range("K76").Value = 100003040
$ws_ALC.Range("M76").Value = -100002725

# ALC row 79
$ws_ALC.Range("H79").Value = 71431890
$ws_ALC.Range("I79").Value = 100003040
$ws_ALC.Range("K79").Value = 100003040
$ws_ALC.Range("M79").Value = -100001948

# ALC row 80
$ws_ALC.Range("H80").Value = 3387.25
$ws_ALC.Range("I80").Value = 1299
$ws_ALC.Range("J80").Value = 3526.4666
$ws_ALC.Range("K80").Value = 3897
$ws_ALC.Range("L80").Value = 10579.3998
$ws_ALC.Range("M80").Value = -2899
$ws_ALC.Range("N80").Value = -12575.3998

# ALC row 83
$ws_ALC.Range("H83").Value = 3387.25
$ws_ALC.Range("I83").Value = 1299
$ws_ALC.Range("J83").Value = 3526.4666
$ws_ALC.Range("K83").Value = 11691
$ws_ALC.Range("L83").Value = 31738.1994
$ws_ALC.Range("M83").Value = -6699
$ws_ALC.Range("N83").Value = -41722.1994

# ALC row 86
$ws_ALC.Range("H86").Value = 62503400
$ws_ALC.Range("J86").Value = 71431910
$ws_ALC.Range("L86").Value = 71431910
$ws_ALC.Range("N86").Value = -71434156

# ALC row 89
$ws_ALC.Range("H89").Value = 62503400
$ws_ALC.Range("J89").Value = 71431910
$ws_ALC.Range("L89").Value = 357159550
$ws_ALC.Range("N89").Value = -357170782

# ALC row 100
$ws_ALC.Range("H100").Value = 4728
$ws_ALC.Range("I100").Value = 4499.75
$ws_ALC.Range("K100").Value = 4499.75
$ws_ALC.Range("M100").Value = -3958.75

# ALC row 112
$ws_ALC.Range("H112").Value = 1169.8667
$ws_ALC.Range("I112").Value = 659.4
$ws_ALC.Range("J112").Value = 1271.96
$ws_ALC.Range("K112").Value = 1978.2
$ws_ALC.Range("L112").Value = 3815.88
$ws_ALC.Range("M112").Value = -870.1999999999998
$ws_ALC.Range("N112").Value = -6031.88

# ALC row 129
$ws_ALC.Range("H129").Value = 1570.2222
$ws_ALC.Range("I129").Value = 1484.9333
$ws_ALC.Range("J129").Value = 1996.6666
$ws_ALC.Range("K129").Value = 4454.7999
$ws_ALC.Range("L129").Value = 5989.9998
$ws_ALC.Range("M129").Value = 545.2001
$ws_ALC.Range("N129").Value = -15989.9998

# ALC row 132
$ws_ALC.Range("H132").Value = 2268.5
$ws_ALC.Range("I132").Value = 2395.5
$ws_ALC.Range("K132").Value = 7186.5
$ws_ALC.Range("M132").Value = -4656.5

# ALC row 138
$ws_ALC.Range("H138").Value = 1412.3658
$ws_ALC.Range("I138").Value = 931
$ws_ALC.Range("K138").Value = 2793
$ws_ALC.Range("M138").Value = 2347

# ALC row 141
$ws_ALC.Range("H141").Value = 8056.3
$ws_ALC.Range("I141").Value = 7223.4287
$ws_ALC.Range("K141").Value = 21670.2861
$ws_ALC.Range("M141").Value = -16490.2861

# ARM row 2
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H2").Value = 707.9
$ws_ARM.Range("I2").Value = 597.375
$ws_ARM.Range("K2").Value = 597.375
$ws_ARM.Range("M2").Value = -484.375

# ARM row 61
$ws_ARM.Range("H61").Value = 1367.3636
$ws_ARM.Range("I61").Value = 1083.2632
$ws_ARM.Range("K61").Value = 1083.2632
$ws_ARM.Range("M61").Value = -871.2632000000001

# ARM row 94
$ws_ARM.Range("H94").Value = 39750
$ws_ARM.Range("J94").Value = 39750
$ws_ARM.Range("L94").Value = 39750
$ws_ARM.Range("N94").Value = -41552

# ARM row 116
$ws_ARM.Range("H116").Value = 707.9
$ws_ARM.Range("I116").Value = 597.375
$ws_ARM.Range("K116").Value = 597.375
$ws_ARM.Range("M116").Value = 1696.625

# ARM row 118
$ws_ARM.Range("H118").Value = 54552.89
$ws_ARM.Range("J118").Value = 54552.89
$ws_ARM.Range("L118").Value = 54552.89
$ws_ARM.Range("N118").Value = -57866.89

# ARM row 124
$ws_ARM.Range("H124").Value = 20014
$ws_ARM.Range("J124").Value = 20014
$ws_ARM.Range("L124").Value = 20014
$ws_ARM.Range("N124").Value = -29834

# ARM row 132
$ws_ARM.Range("H132").Value = 1664.4517
$ws_ARM.Range("I132").Value = 1301.1052
$ws_ARM.Range("J132").Value = 2239.75
$ws_ARM.Range("K132").Value = 3903.3156
$ws_ARM.Range("L132").Value = 6719.25
$ws_ARM.Range("M132").Value = -1373.3156
$ws_ARM.Range("N132").Value = -11779.25

# ARM row 136
$ws_ARM.Range("H136").Value = 1367.3636
$ws_ARM.Range("I136").Value = 1083.2632
$ws_ARM.Range("K136").Value = 3249.7896
$ws_ARM.Range("M136").Value = -699.7896000000001

# BSM row 3
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H3").Value = 707.9
$ws_BSM.Range("I3").Value = 597.375
$ws_BSM.Range("K3").Value = 597.375
$ws_BSM.Range("M3").Value = -483.375

# BSM row 42
$ws_BSM.Range("H42").Value = 349999
$ws_BSM.Range("J42").Value = 349999
$ws_BSM.Range("L42").Value = 349999
$ws_BSM.Range("N42").Value = -350655

# BSM row 99
$ws_BSM.Range("H99").Value = 47527.09
$ws_BSM.Range("I99").Value = 84601.414
$ws_BSM.Range("J99").Value = 3037.9
$ws_BSM.Range("K99").Value = 84601.414
$ws_BSM.Range("L99").Value = 3037.9
$ws_BSM.Range("M99").Value = -83103.414
$ws_BSM.Range("N99").Value = -6033.9

# BSM row 105
$ws_BSM.Range("H105").Value = 94536
$ws_BSM.Range("I105").Value = 201879.8
$ws_BSM.Range("J105").Value = 5082.8335
$ws_BSM.Range("K105").Value = 201879.8
$ws_BSM.Range("L105").Value = 5082.8335
$ws_BSM.Range("M105").Value = -200132.8
$ws_BSM.Range("N105").Value = -8576.833500000001

# BSM row 107
$ws_BSM.Range("H107").Value = 1706.1666
$ws_BSM.Range("I107").Value = 1542.0667
$ws_BSM.Range("K107").Value = 1542.0667
$ws_BSM.Range("M107").Value = 377.9332999999999

# BSM row 134
$ws_BSM.Range("H134").Value = 14362
$ws_BSM.Range("I134").Value = 14362
$ws_BSM.Range("K134").Value = 43086
$ws_BSM.Range("M134").Value = -40551

# BSM row 138
$ws_BSM.Range("H138").Value = 99999
$ws_BSM.Range("J138").Value = 99999
$ws_BSM.Range("L138").Value = 99999
$ws_BSM.Range("N138").Value = -110279

# CRP row 81
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H81").Value = 50000
$ws_CRP.Range("J81").Value = 50000
$ws_CRP.Range("L81").Value = 50000
$ws_CRP.Range("N81").Value = -51996

# CRP row 84
$ws_CRP.Range("H84").Value = 50000
$ws_CRP.Range("J84").Value = 50000
$ws_CRP.Range("L84").Value = 150000
$ws_CRP.Range("N84").Value = -159984

# CRP row 107
$ws_CRP.Range("H107").Value = 1278.3
$ws_CRP.Range("I107").Value = 1377.5
$ws_CRP.Range("K107").Value = 1377.5
$ws_CRP.Range("M107").Value = 542.5

# CRP row 122
$ws_CRP.Range("H122").Value = 2058.9565
$ws_CRP.Range("I122").Value = 1810.4667
$ws_CRP.Range("J122").Value = 2524.875
$ws_CRP.Range("K122").Value = 5431.4001
$ws_CRP.Range("L122").Value = 7574.625
$ws_CRP.Range("M122").Value = -2981.4001
$ws_CRP.Range("N122").Value = -12474.625

# CRP row 132
$ws_CRP.Range("H132").Value = 2200
$ws_CRP.Range("I132").Value = 2200
$ws_CRP.Range("K132").Value = 6600
$ws_CRP.Range("M132").Value = -4070

# CUL row 37
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H37").Value = 74499.75
$ws_CUL.Range("J37").Value = 74499.75
$ws_CUL.Range("L37").Value = 223499.25
$ws_CUL.Range("N37").Value = -223723.25

# CUL row 64
$ws_CUL.Range("H64").Value = 5029.6665
$ws_CUL.Range("I64").Value = 5029.6665
$ws_CUL.Range("K64").Value = 15088.9995
$ws_CUL.Range("M64").Value = -14818.9995

# CUL row 67
$ws_CUL.Range("H67").Value = 5029.6665
$ws_CUL.Range("I67").Value = 5029.6665
$ws_CUL.Range("K67").Value = 15088.9995
$ws_CUL.Range("M67").Value = -14152.9995

# CUL row 97
$ws_CUL.Range("H97").Value = 194.84616
$ws_CUL.Range("I97").Value = 144.7
$ws_CUL.Range("J97").Value = 362
$ws_CUL.Range("K97").Value = 434.1
$ws_CUL.Range("L97").Value = 1086
$ws_CUL.Range("M97").Value = 61.90000000000003
$ws_CUL.Range("N97").Value = -2078

# CUL row 121
$ws_CUL.Range("H121").Value = 2229.6667
$ws_CUL.Range("J121").Value = 2435.0667
$ws_CUL.Range("L121").Value = 7305.2001
$ws_CUL.Range("N121").Value = -9925.2001

# CUL row 140
$ws_CUL.Range("H140").Value = 1246.0667
$ws_CUL.Range("I140").Value = 995.62964
$ws_CUL.Range("K140").Value = 2986.88892
$ws_CUL.Range("M140").Value = 2193.11108

# GSM row 132
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 4219.4287
$ws_GSM.Range("I132").Value = 3379.9285
$ws_GSM.Range("K132").Value = 10139.7855
$ws_GSM.Range("M132").Value = -7609.7855

# LTW row 7
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 4806.4546
$ws_LTW.Range("I7").Value = 3536.5
$ws_LTW.Range("K7").Value = 3536.5
$ws_LTW.Range("M7").Value = -3424.5

# LTW row 68
$ws_LTW.Range("H68").Value = 6001.5
$ws_LTW.Range("I68").Value = 6001.5
$ws_LTW.Range("K68").Value = 6001.5
$ws_LTW.Range("M68").Value = -5252.5

# LTW row 71
$ws_LTW.Range("H71").Value = 6001.5
$ws_LTW.Range("I71").Value = 6001.5
$ws_LTW.Range("K71").Value = 30007.5
$ws_LTW.Range("M71").Value = -26263.5

# LTW row 126
$ws_LTW.Range("H126").Value = 4806.4546
$ws_LTW.Range("I126").Value = 3536.5
$ws_LTW.Range("K126").Value = 10609.5
$ws_LTW.Range("M126").Value = -8139.5

# LTW row 136
$ws_LTW.Range("H136").Value = 2627.9644
$ws_LTW.Range("I136").Value = 3440.7856
$ws_LTW.Range("K136").Value = 10322.3568
$ws_LTW.Range("M136").Value = -7772.356800000001

# WVR row 34
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H34").Value = 0
$ws_WVR.Range("J34").Value = 0
$ws_WVR.Range("N34").Value = ""

# WVR row 100
$ws_WVR.Range("H100").Value = 4763656.5
$ws_WVR.Range("I100").Value = 5954316
$ws_WVR.Range("K100").Value = 11908632
$ws_WVR.Range("M100").Value = -11908091

# WVR row 102
$ws_WVR.Range("H102").Value = 56666.668
$ws_WVR.Range("J102").Value = 56666.668
$ws_WVR.Range("L102").Value = 56666.668
$ws_WVR.Range("N102").Value = -63156.668

# WVR row 136
$ws_WVR.Range("H136").Value = 1752.0769
$ws_WVR.Range("I136").Value = 1481.4166
$ws_WVR.Range("J136").Value = 5000
$ws_WVR.Range("K136").Value = 4444.2498
$ws_WVR.Range("L136").Value = 15000
$ws_WVR.Range("M136").Value = -1894.2498
$ws_WVR.Range("N136").Value = -20100
